$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.1679769670178726,
    0.1813377507281229,
    0.1945435546769042,
    0.2075954808637904,
    0.220494628425907,
    0.2332420936324628,
    0.2458389698792469,
    0.2582863476830883,
    0.2705853146762803,
    0.2827369556009673,
    0.2947423523034943,
    0.3066025837287191,
    0.3183187259142862,
    0.3298918519848619,
    0.3413230321463304,
    0.352613333679951,
    0.3637638209364745,
    0.3747755553302197,
    0.3856495953331089,
    0.3963869964686624,
    0.4069888113059499,
    0.4174560894535014,
    0.4277898775531736,
    0.4379912192739736,
    0.4480611553058388,
    0.4580007233533722,
    0.4678109581295323,
    0.477492891349278,
    0.4870475517231674,
    0.4964759649509092,
    0.5057791537148676,
    0.5149581376735191,
    0.5240139334548602,
    0.5329475546497673,
    0.5417600118053053,
    0.5504523124179873,
    0.559025460926982,
    0.5674804587072706,
    0.5758183040627503,
    0.5840399922192853,
    0.5921465153177041,
    0.6001388624067423,
    0.6080180194359297,
    0.6157849692484221,
    0.6234406915737762,
    0.6309861630206672,
    0.6384223570695474,
    0.6457502440652469,
    0.6529707912095138,
    0.6600849625534934,
    0.6670937189901465,
    0.6739980182466051,
    0.6807988148764645,
    0.6874970602520118,
    0.6940937025563886,
    0.7005896867756878,
    0.7069859546909842,
    0.7132834448702962,
    0.7194830926604783,
    0.725585830179045,
    0.7315925863059223,
    0.7375042866751275,
    0.743321853666376,
    0.749046206396614,
    0.754678260711475,
    0.7602189291766605,
    0.7656691210692425,
    0.7710297423688874,
    0.7763016957489995,
    0.7814858805677827,
    0.7865831928592201,
    0.7915945253239692,
    0.7965207673201717,
    0.8013628048541765,
    0.8061215205711746,
    0.8107977937457439,
    0.8153925002723039,
    0.8199065126554768,
    0.8243407000003544,
    0.8286959280026706,
    0.8329730589388741,
    0.8371729516561046,
    0.8412964615620667,
    0.8453444406148016,
    0.8493177373123548,
    0.8532171966823372,
    0.857043660271379,
    0.8607979661344729,
    0.8644809488242058,
    0.8680934393798765,
    0.8716362653164976,
    0.875110250613679,
    0.878516215704391,
    0.8818549774636055,
    0.8851273491968118,
    0.8883341406284052,
    0.8914761578899462,
    0.8945542035082872,
    0.8975690763935654,
    0.9005215718270571,
    0.9034124814488926,
    0.9062425932456288,
    0.9090126915376751,
    0.9117235569665711,
    0.9143759664821132,
    0.9169706933293252,
    0.9195085070352722,
    0.9219901733957123,
    0.9244164544615839,
    0.926788108525325,
    0.9291058901070203,
    0.9313705499403728,
    0.9335828349584967,
    0.9357434882795258,
    0.9378532491920353,
    0.9399128531402718,
    0.9419230317091873,
    0.9438845126092723,
    0.9457980196611845,
    0.9476642727801668,
    0.9494839879602516,
    0.9512578772582433,
    0.9529866487774774,
    0.9546710066513476,
    0.9563116510265971,
    0.9579092780463671,
    0.9594645798329974,
    0.9609782444705723,
    0.9624509559872051,
    0.9638833943370545,
    0.9652762353820665,
    0.9666301508734341,
    0.9679458084327665,
    0.9692238715329616,
    0.9704649994787724,
    0.9716698473870589,
    0.9728390661667173,
    0.9739733024982773,
    0.975073198813157,
    0.976139393272567,
    0.9771725197460521,
    0.9781732077896605,
    0.9791420826237294,
    0.9800797651102741,
    0.9809868717299707,
    0.9818640145587172,
    0.9827118012437608,
    0.9835308349793781,
    0.9843217144820922,
    0.9850850339654125,
    0.9858213831140811,
    0.9865313470578081,
    0.98721550634448,
    0.987874436912821,
    0.9885087100644887,
    0.9891188924355845,
    0.9897055459675554,
    0.9902692278774667,
    0.9908104906276208,
    0.9913298818944974,
    0.9918279445369884,
    0.9923052165639004,
    0.9927622311006943,
    0.9931995163554322,
    0.9936175955838972,
    0.9940169870538518,
    0.9943982040083986,
    0.9947617546284027,
    0.995108141993934,
    0.9954378640446865,
    0.9957514135393243,
    0.9960492780137067,
    0.9963319397379362,
    0.9965998756721722,
    0.9968535574211482,
    0.9970934511873267,
    0.9973200177226195,
    0.997533712278597,
    0.9977349845551045,
    0.9979242786471938,
    0.998102032990277,
    0.9982686803033946,
    0.9984246475304868,
    0.9985703557795432,
    0.9987062202594978,
    0.9988326502147208,
    0.998950048856948,
    0.9990588132944723,
    0.9991593344584011,
    0.9992519970257688,
    0.9993371793392672,
    0.9994152533233314,
    0.9994865843962916,
    0.9995515313782655,
    0.9996104463944265,
    0.999663674773242,
    0.9997115549392178,
    0.9997544182996295,
    0.9997925891246465,
    0.9998263844201718,
    0.9998561137926204,
    0.9998820793047408,
    0.9999045753214421,
    0.9999238883444207,
    0.9999402968341758,
    0.999954071017757,
    0.9999654726802927,
    0.9999747549379963,
    0.9999821619899254,
    0.9999879288452791,
    0.9999922810224282,
    0.9999954342150651,
    0.9999975939190716,
    0.9999989550063371,
    0.9999997011859848
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $values[$i]
}

Write-Output "Updated A2:A216 with new values"